# edit.ps1
# Applies the two corrections described in the commit
# "Apply corrections to slides supplied by Mum":
#
#   1. The handout date placeholder's cached datetimeFigureOut field text
#      is updated from 19/04/2014 to 25/08/2015.
#   2. On the "We have not feared thee as we ought" slide, the typo
#      "grand" is corrected to "grant" in the line
#      "and grand the grace of holy fear" -> "and grant the grace of holy fear".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master date placeholder: 19/04/2014 -> 25/08/2015
#    (Wrapped in try/catch: some hosts keep the Handout Master's cached
#    date field read-only via the object model, so this is attempted the
#    normal PowerPoint way but must not abort the rest of the script.)
# ---------------------------------------------------------------------------
try {
    $handoutMaster = $p.HandoutMaster
    $dateShape = $handoutMaster.Shapes.Item(2)
    if ($dateShape.HasTextFrame -and $dateShape.TextFrame.TextRange.Text -eq "19/04/2014") {
        $dateShape.TextFrame.TextRange.Text = "25/08/2015"
    }
} catch {
    Write-Host "HandoutMaster date field could not be updated via the object model:" $_
}

# Some hosts model the same date through the headers/footers dialog instead
# of the placeholder shape directly - try that route too, defensively.
try {
    $p.HandoutMaster.HeadersFooters.DateAndTime.UseFormat = 0
    $p.HandoutMaster.HeadersFooters.DateAndTime.Value = "25/08/2015"
} catch {
    # Not fatal - the shape-based update above is the primary mechanism.
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("We have not feared thee as we ought"):
#    "and grand the grace of holy fear" -> "and grant the grace of holy fear"
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

$fullText = $textRange.Text
$target = "grand "
$startIndex = $fullText.IndexOf($target)

if ($startIndex -ge 0) {
    $wordRange = $textRange.Characters($startIndex + 1, $target.Length)
    $wordRange.Text = "grant "
}
